# Add three new API test cases (TC-058-API-19/20/21) covering SMS sending
# behaviour for the booking endpoint, per commit "Added three more test
# cases to test SMS sending."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---- Row 60: booking succeeds and SMS is sent ----
$ws.Cells.Item(60, 1).Value  = "TC-058-API-19"
$ws.Cells.Item(60, 2).Value  = "US-04"
$ws.Cells.Item(60, 3).Value  = "Create booking succeeds and triggers SMS (smsStatus=Sent)"
$ws.Cells.Item(60, 4).Value  = "API"
$ws.Cells.Item(60, 5).Value  = "High"
$ws.Cells.Item(60, 6).Value  = "High"
$ws.Cells.Item(60, 7).Value  = "Server running. Valid service exists. Date/time slot is available (no existing booking for same date+time). Twilio env vars configured. Recipient phone is verified in Twilio (trial)."
$ws.Cells.Item(60, 8).Value  = "1) In Postman, send POST /api/bookings.\n2) Provide valid payload with verified phone.\n3) Send request.\n4) Verify response + SMS received."
$ws.Cells.Item(60, 9).Value  = "serviceId=<VALID_SERVICE_ID>; date=2026-01-20; time=10:00; customerName=Robert Norwood; phone=<VERIFIED_E164_PHONE>; email=test@example.com; notes=Success path"
$ws.Cells.Item(60, 10).Value = "Returns 201. Response success=true and booking is created. smsStatus=Sent. smsError empty. SMS confirmation received on phone."
$ws.Cells.Item(60, 13).Value = "Robert Norwood"
$ws.Cells.Item(60, 14).Value = 46041
$ws.Cells.Item(60, 15).Value = "Manual API test via Postman. Confirms booking create + Twilio integration works end-to-end."

# ---- Row 61: booking rejects duplicate date/time slot ----
$ws.Cells.Item(61, 1).Value  = "TC-058-API-20"
$ws.Cells.Item(61, 2).Value  = "US-04"
$ws.Cells.Item(61, 3).Value  = "Create booking rejects duplicate date/time slot(409)."
$ws.Cells.Item(61, 4).Value  = "API"
$ws.Cells.Item(61, 5).Value  = "High"
$ws.Cells.Item(61, 6).Value  = "Medium"
$ws.Cells.Item(61, 7).Value  = "Server running. Valid service exists. A booking already exists for date=2026-01-20 and time=10:00 (same slot)."
$ws.Cells.Item(61, 8).Value  = "1) Ensure an existing booking already uses date=2026-01-20 and time=10:00.\n2) In Postman, send POST /api/bookings with the same date/time.\n3) Send request.\n4) Verify conflict response and no new booking created."
$ws.Cells.Item(61, 9).Value  = "serviceId=<VALID_SERVICE_ID>; date=2026-01-20; time=10:00; customerName=Jane Doe; phone=<VERIFIED_E164_PHONE>; email=dup@test.com; notes=Duplicate slot test"
$ws.Cells.Item(61, 10).Value = "Returns 409 with message `"Slot already booked.`" No new booking created. No SMS sent for the failed request."
$ws.Cells.Item(61, 13).Value = "Robert Norwood"
$ws.Cells.Item(61, 14).Value = 46041
$ws.Cells.Item(61, 15).Value = "Manual API test via Postman. Validates unique index enforcement on {date, time}."

# ---- Row 62: booking still succeeds when SMS fails ----
$ws.Cells.Item(62, 1).Value  = "TC-058-API-21"
$ws.Cells.Item(62, 2).Value  = "US-05"
$ws.Cells.Item(62, 3).Value  = "Booking succeeds even if SMS fails (smsStatus=Failed, smsError set)"
$ws.Cells.Item(62, 4).Value  = "API"
$ws.Cells.Item(62, 5).Value  = "High"
$ws.Cells.Item(62, 6).Value  = "High"
$ws.Cells.Item(62, 7).Value  = "Server running. Valid service exists. Date/time slot is available. Force SMS failure by using an unverified phone number OR temporarily invalid TWILIO credentials OR remove TWILIO_MESSAGING_SERVICE_SID."
$ws.Cells.Item(62, 8).Value  = "1) Intentionally force SMS failure (unverified phone or invalid Twilio env).\n2) In Postman, send POST /api/bookings.\n3) Send request.\n4) Verify booking is still created and smsStatus reflects failure."
$ws.Cells.Item(62, 9).Value  = "serviceId=<VALID_SERVICE_ID>; date=2026-01-21; time=11:00; customerName=SMS Fail Test; phone=<UNVERIFIED_OR_INVALID_E164_PHONE>; email=smsfail@test.com; notes=Force SMS failure"
$ws.Cells.Item(62, 10).Value = "Returns 201. Booking is created successfully (status remains Confirmed). smsStatus=Failed. smsError contains Twilio error message (or fallback). User does not need SMS success for booking creation."
$ws.Cells.Item(62, 13).Value = "Robert Norwood"
$ws.Cells.Item(62, 14).Value = 46041
$ws.Cells.Item(62, 15).Value = "Manual API test via Postman. Validates non-blocking SMS requirement (US-05)."

# ---- View state: zoom in to 120% and leave selection on the last new row ----
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("G62").Select()
